$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.178.80"
$ws.Range("E2").Value = "  -0.31%  "

$ws.Range("D3").Value = "3.803.07"
$ws.Range("E3").Value = "  +1.02%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.48%  "

$ws.Range("D7").Value = "3.801.26"
$ws.Range("E7").Value = "  +0.98%  "

$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("E9").Value = "  -0.60%  "

$ws.Range("E10").Value = "  +0.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.44%  "

$ws.Range("E12").Value = "  -0.70%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.19%  "

$ws.Range("E14").Value = "  -2.16%  "

$ws.Range("D15").Value = "4.438.58"
$ws.Range("E15").Value = "  +1.10%  "

$ws.Range("D16").Value = "3.810.28"
$ws.Range("E16").Value = "  +1.36%  "

$ws.Range("D17").Value = "69.272.51"
$ws.Range("E17").Value = "  -0.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.45%  "

$ws.Range("E19").Value = "  -0.24%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "487.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.78%  "

$ws.Range("E23").Value = "  -1.08%  "

$ws.Range("E24").Value = "  -1.99%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.74"
$ws.Range("D25").Style = "Normal"

$ws.Range("E26").Value = "  -3.51%  "

$ws.Range("E27").Value = "  -1.23%  "

$ws.Range("E28").Value = "  -2.42%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("E30").Value = "  -0.71%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.04%  "

$ws.Range("E32").Value = "  -5.24%  "

$ws.Range("D33").Value = "3.948.90"
$ws.Range("E33").Value = "  +1.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.11%  "

$ws.Range("D35").Value = "3.749.17"
$ws.Range("E35").Value = "  +1.41%  "

$ws.Range("E36").Value = "  -2.00%  "

$ws.Range("E37").Value = "  +5.71%  "

$ws.Range("E38").Value = "  -0.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.15%  "

$ws.Range("E41").Value = "  -0.88%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.31%  "

$ws.Range("E43").Value = "  +0.82%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "423.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.26%  "

$ws.Range("E47").Value = "  -1.39%  "

$ws.Range("D48").Value = "2.823.93"
$ws.Range("E48").Value = "  +0.30%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "39.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.95%  "

$ws.Range("E51").Value = "  -1.53%  "
